$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.675.22"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.851.63"
$ws.Range("E3").Value = "  -3.69%  "
$ws.Range("E4").Value = "  -0.95%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.84"
$ws.Range("E5").Value = "  +2.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4646"
$ws.Range("E7").Value = "  -3.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3914"
$ws.Range("E8").Value = "  -3.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.51"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07929"
$ws.Range("E10").Value = "  -3.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9858"
$ws.Range("E11").Value = "  -2.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.32"
$ws.Range("E12").Value = "  -5.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.929.47"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.851"
$ws.Range("E14").Value = "  -3.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.003"
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06850"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("E17").Value = "  -4.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001008"
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("E20").Value = "  -2.81%  "
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.700.94"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.396"
$ws.Range("E23").Value = "  -4.99%  "
$ws.Range("E24").Value = "  -5.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.182.65"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.127"
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.15"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.48"
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.098"
$ws.Range("E29").Value = "  -5.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.025"
$ws.Range("E30").Value = "  -3.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.73"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9794"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09417"
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.371"
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.481"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.355"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06160"
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02206"
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.164"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.000"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5725"
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.634"
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.17"
$ws.Range("E43").Value = "  -5.79%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1801"
$ws.Range("E44").Value = "  -2.69%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.394"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.247"
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.87"
$ws.Range("E47").Value = "  -4.66%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5409"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07160"
$ws.Range("E49").Value = "  -4.35%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.911"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.90"
$ws.Range("E51").Value = "  -4.46%  "
